# Edit: split the "با پیش روی در بازی نوع، ..." paragraph into several
# runs, inserting "؛" right after "با پیش روی در بازی" and "با " right
# before "زدن سکه ها ...", while leaving the existing "_GoBack" bookmark
# positioned between the new "با " run and the "زدن ..." run (matching
# the target OOXML, where the bookmark sits directly after that run).

$d = $word.ActiveDocument

$anchor = "با پیش روی در بازی نوع"
$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text.Contains($anchor)) {
        $target = $para
        break
    }
}
if ($target -eq $null) {
    throw "Could not locate target paragraph"
}

$start = $target.Range.Start

# Offsets (relative to paragraph start) computed from the original text:
#   "با پیش روی در بازی" | " نوع، ..."
#   -> insert "؛" right after "بازی" (i.e. right before the space + نوع)
$prefix = "با پیش روی در بازی"
$semiPos = $start + $prefix.Length

# Remove the existing bookmark up front so it doesn't anchor itself to
# "end of paragraph text" while we are splicing new text in - we'll
# recreate it afterwards at the exact spot the diff wants.
$hadBookmark = $d.Bookmarks.Exists("_GoBack")
if ($hadBookmark) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Insert "؛" right after "بازی"
$d.Range($semiPos, $semiPos).InsertAfter([string]([char]0x061B))

# Recompute, then insert "با " right before "زدن"
$fullText = $target.Range.Text
$zadanIdx = $fullText.IndexOf("زدن")
$baPos = $start + $zadanIdx
$d.Range($baPos, $baPos).InsertAfter("با ")

# Compute the run-split boundaries (relative offsets, in the now-updated text):
#   0 .. len(prefix)                      -> "با پیش روی در بازی"
#   .. +1                                 -> "؛"
#   .. up to (new) زدن-با insertion point -> " نوع، ... می تواند "
#   .. +3                                 -> "با "
#   .. end                                -> "زدن ... بپردازد."
$b1 = $prefix.Length               ; # end of "با پیش روی در بازی"
$b2 = $b1 + 1                      ; # end of "؛"
$b3 = $baPos - $start              ; # end of the middle chunk (start of "با ")
$b4 = $b3 + 3                      ; # end of "با "

$paraEnd = $target.Range.End        # includes the trailing paragraph mark

function Split-RunAt($relPos) {
    # A range that spans from the split point to the end of the paragraph
    # text (but excluding the paragraph mark) forces the engine to break
    # the run boundary only at $relPos, since the far end already sits on
    # an existing boundary (the paragraph end).
    $r = $d.Range($start + $relPos, $paraEnd - 1)
    $r.Bold = $true
    $r.Bold = $false
}

# Apply right-to-left so earlier offsets stay valid (none of these calls
# change the text length, but keeping the order consistent is simplest).
Split-RunAt $b4
Split-RunAt $b3
Split-RunAt $b2
Split-RunAt $b1

# Recreate the bookmark exactly between the new "با " run and the
# "زدن ..." run.
if ($hadBookmark) {
    $bmRange = $d.Range($start + $b4, $start + $b4)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

Write-Output "OK"
